$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[-, -, 'MCT-3A-Lab. de eletroeletrônica']"
$ws.Range("F2").Value = "[-, 'MCT-2A-Sistemas digitais', 'MCT-2A-Sistemas digitais']"

$ws.Range("C4").Value = "-"

$ws.Range("C6").Value = "-"

$ws.Range("C7").Value = "-"

$ws.Range("E8").Value = "[-, -, 'MCT-3A-Lab. de eletroeletrônica']"

$ws.Range("F19").Value = "[-, 'ELM-2NA-Lab. Circuitos Elétricos']"

$ws.Range("D21").Value = "[-, 'ELM-2NA-Lab. Circuitos Elétricos']"
